# Minor text optimizations on the "Einstieg Accessibility" (slide 16) and
# "Landmarks (Orientierungspunkte)" (slide 17) slides.

$p = $ppt.ActivePresentation

# --- Slide 16: "Content Placeholder 2" --------------------------------
$slide16 = $p.Slides.Item(16)
$body16  = $slide16.Shapes.Item(2)
$tf16    = $body16.TextFrame
$tr16    = $tf16.TextRange

# Paragraph 5 = "Focusable Elements (Fokussierbare Elemente)"
$focusPara = $tr16.Paragraphs(5)
$focusPara.Runs(1).Text = "Focusable / Clickable Elements "
$focusPara.Runs(2).Text = "(Fokussierbare / Klickbare Elemente)"

# The placeholder text grew a little, so PowerPoint re-fits it to the
# shape: the previous line-spacing-only reduction becomes a font-scale
# reduction (92.5%) instead.
$tf16.AutofitFontScale = 92.5
$tf16.AutofitLineSpaceReduction = 0

# --- Slide 17: "Inhaltsplatzhalter 2" ----------------------------------
$slide17 = $p.Slides.Item(17)
$body17  = $slide17.Shapes.Item(2)
$tr17    = $body17.TextFrame.TextRange

# Paragraph 2 = "Screen-Reader nutzen die Landmarks fuer die Navigation"
$navPara = $tr17.Paragraphs(2)
$navPara.Runs(2).Text = "Navigation (VO Rotor Landmarks)"
